# Add 2022-Q3 data
# 1) Update the "总计" (Total) summary sheet: insert a new top data row for
#    2022-Q3 and push the existing quarters down by one row.
# 2) Insert a brand-new "2022-Q3" worksheet (holding the per-fund holdings
#    table) positioned right after "总计" and before "2022-Q2".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: "总计" sheet - shift the quarterly summary rows down and fill
# in the new 2022-Q3 figures at the top of the data.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Grow the data block by one row (row 9) using the same cell formatting
# (style) as the existing last row's index cell.
$total.Range("A8").Copy($total.Range("A9"))

$summaryRows = @(
    @(0, "2022-Q3", 7, 0.24),
    @(1, "2022-Q2", 2, 0.27),
    @(2, "2022-Q1", 4, 8.800000000000001),
    @(3, "2021-Q4", 2, 0.59),
    @(4, "2021-Q3", 10, 3.51),
    @(5, "2021-Q2", 6, 1.82),
    @(6, "2021-Q1", 3, 5.14),
    @(7, "2020-Q4", 1, 0.18)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = $i + 2
    $vals = $summaryRows[$i]
    $total.Cells.Item($r, 1).Value = $vals[0]
    $total.Cells.Item($r, 2).Value = $vals[1]
    $total.Cells.Item($r, 3).Value = $vals[2]
    $total.Cells.Item($r, 4).Value = $vals[3]
}

# ---------------------------------------------------------------------
# Step 2: Create the new "2022-Q3" sheet. The quarterly sheets all share
# the same layout, so clone the existing "2022-Q2" sheet (it carries the
# right column headers/styles) and place the clone right before it, then
# rename it and replace its fund-holding data with the 2022-Q3 figures.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Extend the index column (A) formatting down to the 7 fund rows needed.
$q3.Range("A2").Copy($q3.Range("A3:A8"))

$funds = @(
    @("233009", "大摩多因子精选策略混合", "6.50", "83.44", "0.88", "0.0572", 10),
    @("008115", "天弘中证红利低波动100指数C", "2.44", "94.56", "2.03", "0.0495", 4),
    @("159617", "华夏中证智选500价值稳健策略ETF", "2.93", "97.05", "1.44", "0.0422", 3),
    @("008114", "天弘中证红利低波动100指数A", "1.89", "94.56", "2.03", "0.0384", 4),
    @("515100", "景顺长城中证红利低波动100ETF", "1.62", "98.63", "2.13", "0.0345", 4),
    @("009658", "汇丰晋信中小盘低波动策略股票A", "0.85", "90.14", "1.97", "0.0167", 3),
    @("009775", "汇丰晋信中小盘低波动策略股票C", "0.04", "90.14", "1.97", "0.0008", 3)
)

for ($i = 0; $i -lt $funds.Length; $i++) {
    $r = $i + 2
    $fund = $funds[$i]
    $q3.Cells.Item($r, 1).Value = $i
    $q3.Cells.Item($r, 2).Value = "'" + $fund[0]
    $q3.Cells.Item($r, 3).Value = $fund[1]
    $q3.Cells.Item($r, 4).Value = "'" + $fund[2]
    $q3.Cells.Item($r, 5).Value = "'" + $fund[3]
    $q3.Cells.Item($r, 6).Value = "'" + $fund[4]
    $q3.Cells.Item($r, 7).Value = "'" + $fund[5]
    $q3.Cells.Item($r, 8).Value = $fund[6]
}
